# MMD json data tweaks (appendix_changes.xlsx)
# ---------------------------------------------------------------------
# Applies the appendix edits to Sheet1:
#  - fixes the stray newline in the Syria "Title:" snippet
#  - fills in two missing "Added text?" snippets (pages 9 and 30)
#  - clears the red "still open" flag on pages 9 and 30 now that they
#    have snippets
#  - splits the page-73 Syria entry into a merged two-row block (a
#    "Text:" / "Text-end:" pair) and answers the old "should 74, 75, 76
#    be title?" question directly in the sheet for pages 74-76
# ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Fix the Syria "Title:" text in D16 -- drop the stray leading newline.
$ws.Range("D16").Value = "Title: The humanitarian fallout of the conflict in Syria reaches new proportions as the number of estimated refugees reaches one million. Who's Helping?"

# 2. Fill in the missing "Added text?" snippets for page 9 (row 7) and
#    page 30 (row 18).
$ws.Range("D7").Value = "Text: The Pew Forum survey included several questions designed to probe the kinds of requests that inmates make for accommodation of their religious beliefs and practices, as well as the frequency with which they are granted."
$ws.Range("D18").Value = "Text: The idea that hard work leads to material success is no longer, if it ever was, a uniquely Western value."

# 3. Those rows are no longer "still open" -- clear the red flag font and
#    give them the normal yellow-fill / centered page-number look used
#    elsewhere in column A.
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. Insert a row after row 21 (page 73) so the Syria entry becomes a
#    merged two-row block, pushing old rows 22/23/25 down to 23/24/26.
$ws.Rows.Item(22).Insert()

# Page 73 also loses its red "still open" flag now that it has text --
# give the new merged block the same look as the rest of column A.
$ws.Range("A2").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A21:A22").Merge()
$ws.Range("B21:B22").Merge()
$ws.Range("C21:C22").Merge()
$ws.Range("A21:A22").VerticalAlignment = -4108
$ws.Range("B21:B22").HorizontalAlignment = -4108
$ws.Range("C21:C22").HorizontalAlignment = -4108

$ws.Range("D21").Value = "Text: Two years after nationwide protests sparked upheaval in Syria, the ensuing refugee crisis has reached one million people. "
$ws.Range("D22").Value = "Text-end: The British charity Save the Children claims that many of these children have been separated from one or both of their parents."

# 5. Rows 23/24 (pages 75 and 76, shifted down from 22/23) turn out to be
#    about the same Syria title -- fill in column D for both, matching
#    the wrap formatting used for the other "Title:" entries.
$ws.Range("D16").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D23").Value = "Title: The humanitarian fallout of the conflict in Syria reaches new proportions as the number of estimated refugees reaches one million. Who's Helping?"
$ws.Range("D24").Value = "Title: The humanitarian fallout of the conflict in Syria reaches new proportions as the number of estimated refugees reaches one million. Who's Helping?"

# 6. The old "Should 74, 75, 76 be title?" placeholder question (now on
#    row 27 after the insert) is answered by rows 23/24 above, so replace
#    it with the "NOTES: / 3 has a new reference" note that used to sit
#    on row 25 (now row 26).
$ws.Range("D26").Value = "3 has a new reference"
$ws.Range("D27").ClearContents()

# 7. Restore the selection to where editing left off.
$ws.Range("D31").Select()
